$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above row 95; existing rows 95-102 shift down to 96-103.
$ws.Rows("95").Insert()

# Populate the newly inserted row 95 with the new weekly record.
$ws.Range("A95").Value = 9
$ws.Range("B95").Value = "Vega Central Mapocho de Santiago"
$ws.Range("C95").Value = "Metropolitana"
$ws.Range("D95").Value = 44769
$ws.Range("E95").Value = 13
$ws.Range("F95").Value = 100112005
$ws.Range("G95").Value = "Puerro"
$ws.Range("H95").Value = "Sin especificar"
$ws.Range("I95").Value = "Primera"
$ws.Range("J95").Value = 70
$ws.Range("K95").Value = 7000
$ws.Range("L95").Value = 8000
$ws.Range("M95").Value = 7500
$ws.Range("N95").Value = "$/paquete 20 unidades"
$ws.Range("O95").Value = "Provincia de Chacabuco"
$ws.Range("P95").Value = 375
$ws.Range("Q95").Value = 20
$ws.Range("R95").Value = "Hortaliza"
